$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 2. Edit old git entry with new category name (column A re-labelled) ---
$ws.Range('A6').Value = 'git-branch'
$ws.Range('A7').Value = 'git-branch'
$ws.Range('A8').Value = 'git-branch'
$ws.Range('A9').Value = 'git-branch'
$ws.Range('A10').Value = 'git-push/pull'
$ws.Range('A11').Value = 'git-log'
$ws.Range('A12').Value = 'git-merge'
$ws.Range('A13').Value = 'git-diff'
$ws.Range('A14').Value = 'git-diff'
$ws.Range('A24').Value = 'git-stash'
$ws.Range('A25').Value = 'git-log'
$ws.Range('A26').Value = 'git-log'

# --- 2. Add git entry on 2016 trend, worktree and rebase (new rows 33-35) ---
$ws.Range('A33').Value = '2016 Overview'
$ws.Range('B33').Value = 'wrap up the 2016 new feature and changes'
$ws.Range('C33').Value = 'According to linux.cn post 2016 Git 新視界, there are below big updates:
1. git worktree improved (more options, auto-rebase...)
2. git rebase improved (more options)
3. git lfs improved (Large file system)
4. git diff improved (new algorithm to avoid miss marking of same wording old and new line)
5. git stash improved
6. git submodule is a pain in the ass
7. Lastest version rolling to v2.11.0'

$ws.Range('A34').Value = 'git-rebase'
$ws.Range('B34').Value = 'Basic intro'
$ws.Range('C34').Value = 'Say A>B>C>D>E is master, A>B>C>G>H is fix, For some reasons, we want to refactor the repo to a linear order (from a so-called "disarray" order), then we perfrom below command to rebase the fix to master base:
$ (at fix branch) git rebase master
then the repo becomes: A>B>C>D>E>G''>H'' (the base of fix changed from C to E)'

$ws.Range('A35').Value = 'git-worktree'
$ws.Range('B35').Value = 'Basic intro'
$ws.Range('C35').Value = 'Git introduced the worktree feature not too long ago (as of version 2.5, released July 2015). A great usage scenario can be found here: https://spin.atomicobject.com/2016/06/26/parallelize-development-git-worktrees/
Set up worktree with below command:
$ git worktree add ../new-worktree-dir some-existing-branch'

# Row heights for the new rows (match source row-height metadata)
$ws.Rows.Item(33).RowHeight = 102.75
$ws.Rows.Item(34).RowHeight = 64.5
$ws.Rows.Item(35).RowHeight = 64.5

# Restore the active selection to where the edits left off
$ws.Range('C36').Select()
